# Set the "Status Processo" value for the first data row (row 2) from
# "FATURADO" to "PENDENTE" - separating out pending receipts from
# already-invoiced ones, per the commit message about commissions by receipt.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "PENDENTE"
